$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 34817
$ws.Cells.Item(2, 4).Value = 50427466
$ws.Cells.Item(3, 3).Value = 85364
$ws.Cells.Item(3, 4).Value = 125296936
$ws.Cells.Item(4, 3).Value = 29286
$ws.Cells.Item(4, 4).Value = 43413901
$ws.Cells.Item(5, 3).Value = 8059
$ws.Cells.Item(5, 4).Value = 11983018
$ws.Cells.Item(6, 3).Value = 1726
$ws.Cells.Item(6, 4).Value = 2566096
$ws.Cells.Item(7, 3).Value = 126
$ws.Cells.Item(7, 4).Value = 184093
$ws.Cells.Item(11, 3).Value = 38243
$ws.Cells.Item(11, 4).Value = 52032936
$ws.Cells.Item(12, 3).Value = 8965
$ws.Cells.Item(12, 4).Value = 12982130
$ws.Cells.Item(13, 3).Value = 24584
$ws.Cells.Item(13, 4).Value = 36082440
$ws.Cells.Item(14, 3).Value = 7843
$ws.Cells.Item(14, 4).Value = 11650450
$ws.Cells.Item(15, 3).Value = 1988
$ws.Cells.Item(15, 4).Value = 2957476
$ws.Cells.Item(19, 3).Value = 9496
$ws.Cells.Item(19, 4).Value = 12625492
$ws.Cells.Item(20, 3).Value = 12524
$ws.Cells.Item(20, 4).Value = 18098022
$ws.Cells.Item(21, 3).Value = 30007
$ws.Cells.Item(21, 4).Value = 44091941
$ws.Cells.Item(22, 3).Value = 9741
$ws.Cells.Item(22, 4).Value = 14487810
$ws.Cells.Item(23, 3).Value = 2467
$ws.Cells.Item(23, 4).Value = 3670763
$ws.Cells.Item(24, 3).Value = 432
$ws.Cells.Item(24, 4).Value = 642345
$ws.Cells.Item(26, 3).Value = 10939
$ws.Cells.Item(26, 4).Value = 14682728
$ws.Cells.Item(27, 3).Value = 7127
$ws.Cells.Item(27, 4).Value = 10332123
$ws.Cells.Item(28, 3).Value = 21248
$ws.Cells.Item(28, 4).Value = 31218892
$ws.Cells.Item(29, 3).Value = 7389
$ws.Cells.Item(29, 4).Value = 10999151
$ws.Cells.Item(30, 3).Value = 1818
$ws.Cells.Item(30, 4).Value = 2715486
$ws.Cells.Item(31, 3).Value = 297
$ws.Cells.Item(31, 4).Value = 442915
$ws.Cells.Item(33, 3).Value = 7737
$ws.Cells.Item(33, 4).Value = 10259641
$ws.Cells.Item(34, 3).Value = 2777
$ws.Cells.Item(34, 4).Value = 4006336
$ws.Cells.Item(35, 3).Value = 6996
$ws.Cells.Item(35, 4).Value = 10224347
$ws.Cells.Item(36, 3).Value = 2836
$ws.Cells.Item(36, 4).Value = 4199023
$ws.Cells.Item(38, 3).Value = 131
$ws.Cells.Item(38, 4).Value = 195804
$ws.Cells.Item(40, 3).Value = 2158
$ws.Cells.Item(40, 4).Value = 2909666
$ws.Cells.Item(41, 3).Value = 16049
$ws.Cells.Item(41, 4).Value = 23228873
$ws.Cells.Item(42, 3).Value = 48219
$ws.Cells.Item(42, 4).Value = 70748433
$ws.Cells.Item(43, 3).Value = 18059
$ws.Cells.Item(43, 4).Value = 26831123
$ws.Cells.Item(44, 3).Value = 5242
$ws.Cells.Item(44, 4).Value = 7813608
$ws.Cells.Item(45, 3).Value = 1039
$ws.Cells.Item(45, 4).Value = 1549719
$ws.Cells.Item(49, 3).Value = 15499
$ws.Cells.Item(49, 4).Value = 20707096
$ws.Cells.Item(51, 3).Value = 5962
$ws.Cells.Item(51, 4).Value = 8778326
$ws.Cells.Item(52, 3).Value = 2086
$ws.Cells.Item(52, 4).Value = 3115750
$ws.Cells.Item(56, 3).Value = 5452
$ws.Cells.Item(56, 4).Value = 7532638
$ws.Cells.Item(57, 3).Value = 685
$ws.Cells.Item(57, 4).Value = 1003840
$ws.Cells.Item(58, 3).Value = 1740
$ws.Cells.Item(58, 4).Value = 2578414
$ws.Cells.Item(59, 3).Value = 695
$ws.Cells.Item(59, 4).Value = 1036719
$ws.Cells.Item(60, 3).Value = 238
$ws.Cells.Item(60, 4).Value = 356758
$ws.Cells.Item(61, 3).Value = 49
$ws.Cells.Item(61, 4).Value = 73500
$ws.Cells.Item(63, 3).Value = 1007
$ws.Cells.Item(63, 4).Value = 1428265
$ws.Cells.Item(64, 3).Value = 14295
$ws.Cells.Item(64, 4).Value = 20666944
$ws.Cells.Item(65, 3).Value = 42345
$ws.Cells.Item(65, 4).Value = 62016327
$ws.Cells.Item(66, 3).Value = 14946
$ws.Cells.Item(66, 4).Value = 22224375
$ws.Cells.Item(67, 3).Value = 4294
$ws.Cells.Item(67, 4).Value = 6397793
$ws.Cells.Item(68, 3).Value = 827
$ws.Cells.Item(68, 4).Value = 1231273
$ws.Cells.Item(71, 3).Value = 14179
$ws.Cells.Item(71, 4).Value = 18774582
$ws.Cells.Item(72, 3).Value = 45543
$ws.Cells.Item(72, 4).Value = 66310164
$ws.Cells.Item(73, 3).Value = 132252
$ws.Cells.Item(73, 4).Value = 194974756
$ws.Cells.Item(74, 3).Value = 58197
$ws.Cells.Item(74, 4).Value = 86751683
$ws.Cells.Item(75, 3).Value = 18460
$ws.Cells.Item(75, 4).Value = 27587072
$ws.Cells.Item(76, 3).Value = 4111
$ws.Cells.Item(76, 4).Value = 6143120
$ws.Cells.Item(83, 3).Value = 45014
$ws.Cells.Item(83, 4).Value = 61594241
$ws.Cells.Item(84, 3).Value = 4078
$ws.Cells.Item(84, 4).Value = 5915259
$ws.Cells.Item(85, 3).Value = 10578
$ws.Cells.Item(85, 4).Value = 15549460
$ws.Cells.Item(86, 3).Value = 3623
$ws.Cells.Item(86, 4).Value = 5400581
$ws.Cells.Item(88, 3).Value = 259
$ws.Cells.Item(88, 4).Value = 386612
$ws.Cells.Item(91, 3).Value = 4747
$ws.Cells.Item(91, 4).Value = 6405537
$ws.Cells.Item(92, 3).Value = 1384
$ws.Cells.Item(92, 4).Value = 2000490
$ws.Cells.Item(93, 3).Value = 4580
$ws.Cells.Item(93, 4).Value = 6749465
$ws.Cells.Item(94, 3).Value = 1776
$ws.Cells.Item(94, 4).Value = 2648616
$ws.Cells.Item(95, 3).Value = 624
$ws.Cells.Item(95, 4).Value = 935141
$ws.Cells.Item(96, 3).Value = 155
$ws.Cells.Item(96, 4).Value = 232069
$ws.Cells.Item(99, 3).Value = 3035
$ws.Cells.Item(99, 4).Value = 4030233
$ws.Cells.Item(100, 3).Value = 514
$ws.Cells.Item(100, 4).Value = 766464
$ws.Cells.Item(101, 3).Value = 264
$ws.Cells.Item(101, 4).Value = 394265
$ws.Cells.Item(105, 3).Value = 10038
$ws.Cells.Item(105, 4).Value = 14586706
$ws.Cells.Item(106, 3).Value = 27745
$ws.Cells.Item(106, 4).Value = 40792721
$ws.Cells.Item(107, 3).Value = 9304
$ws.Cells.Item(107, 4).Value = 13839476
$ws.Cells.Item(108, 3).Value = 2529
$ws.Cells.Item(108, 4).Value = 3770910
$ws.Cells.Item(109, 3).Value = 431
$ws.Cells.Item(109, 4).Value = 643982
$ws.Cells.Item(112, 3).Value = 9145
$ws.Cells.Item(112, 4).Value = 12123046
$ws.Cells.Item(113, 3).Value = 28190
$ws.Cells.Item(113, 4).Value = 40701511
$ws.Cells.Item(114, 3).Value = 62502
$ws.Cells.Item(114, 4).Value = 91554367
$ws.Cells.Item(115, 3).Value = 20262
$ws.Cells.Item(115, 4).Value = 30128838
$ws.Cells.Item(116, 3).Value = 5666
$ws.Cells.Item(116, 4).Value = 8446361
$ws.Cells.Item(117, 3).Value = 1013
$ws.Cells.Item(117, 4).Value = 1515493
$ws.Cells.Item(121, 3).Value = 24030
$ws.Cells.Item(121, 4).Value = 32193758
$ws.Cells.Item(122, 3).Value = 32865
$ws.Cells.Item(122, 4).Value = 47490318
$ws.Cells.Item(123, 3).Value = 71670
$ws.Cells.Item(123, 4).Value = 104916676
$ws.Cells.Item(124, 3).Value = 22435
$ws.Cells.Item(124, 4).Value = 33311954
$ws.Cells.Item(125, 3).Value = 5925
$ws.Cells.Item(125, 4).Value = 8812044
$ws.Cells.Item(126, 3).Value = 1075
$ws.Cells.Item(126, 4).Value = 1602769
$ws.Cells.Item(128, 3).Value = 15
$ws.Cells.Item(128, 4).Value = 22500
$ws.Cells.Item(130, 3).Value = 29136
$ws.Cells.Item(130, 4).Value = 38824537
$ws.Cells.Item(131, 3).Value = 12313
$ws.Cells.Item(131, 4).Value = 17835177
$ws.Cells.Item(132, 3).Value = 30655
$ws.Cells.Item(132, 4).Value = 45059165
$ws.Cells.Item(133, 3).Value = 10930
$ws.Cells.Item(133, 4).Value = 16241690
$ws.Cells.Item(134, 3).Value = 2759
$ws.Cells.Item(134, 4).Value = 4114791
$ws.Cells.Item(135, 3).Value = 439
$ws.Cells.Item(135, 4).Value = 652490
$ws.Cells.Item(138, 3).Value = 10128
$ws.Cells.Item(138, 4).Value = 13568951
$ws.Cells.Item(139, 3).Value = 32105
$ws.Cells.Item(139, 4).Value = 46400265
$ws.Cells.Item(140, 3).Value = 75887
$ws.Cells.Item(140, 4).Value = 111256616
$ws.Cells.Item(141, 3).Value = 22908
$ws.Cells.Item(141, 4).Value = 34066153
$ws.Cells.Item(142, 3).Value = 5900
$ws.Cells.Item(142, 4).Value = 8807348
$ws.Cells.Item(143, 3).Value = 1275
$ws.Cells.Item(143, 4).Value = 1900006
$ws.Cells.Item(146, 3).Value = 27010
$ws.Cells.Item(146, 4).Value = 36607813
